$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 12.2218579614066
$ws.Range("F2").Value = 1.881011407927113
$ws.Range("G2").Value = 0.002583783132887318
$ws.Range("E3").Value = 12.28454731732566
$ws.Range("F3").Value = 1.831409206878172
$ws.Range("G3").Value = -0.002532245904274699
$ws.Range("E4").Value = 12.30414225166378
$ws.Range("F4").Value = 1.761698548583134
$ws.Range("G4").Value = -0.004131373086001311
$ws.Range("E5").Value = 12.41600002994138
$ws.Range("F5").Value = 1.743177915011001
$ws.Range("G5").Value = -0.01325999840541781
$ws.Range("E6").Value = 11.93222199747857
$ws.Range("F6").Value = 1.914494212867257
$ws.Range("G6").Value = 0.02622074637710015
$ws.Range("E7").Value = 12.20646435874626
$ws.Range("F7").Value = 1.836603404368611
$ws.Range("G7").Value = 0.003840042940350008
$ws.Range("E8").Value = 11.89655116997449
$ws.Range("F8").Value = 1.768092597376504
$ws.Range("G8").Value = 0.02913181455789426
$ws.Range("E9").Value = 12.31234066573045
$ws.Range("F9").Value = 1.736308884669101
$ws.Range("G9").Value = -0.004800439210686402
$ws.Range("E10").Value = 11.80416136678433
$ws.Range("F10").Value = 1.949901181618349
$ws.Range("G10").Value = 0.03667167374018976
$ws.Range("E11").Value = 12.13819087238466
$ws.Range("F11").Value = 1.839736334788845
$ws.Range("G11").Value = 0.009411788471524463
$ws.Range("E12").Value = 11.67651281041032
$ws.Range("F12").Value = 1.766757441934289
$ws.Range("G12").Value = 0.04708897204206464
$ws.Range("E13").Value = 12.22302506305061
$ws.Range("F13").Value = 1.736163274384311
$ws.Range("G13").Value = 0.002488536885536852
$ws.Range("E14").Value = 11.82060991876815
$ws.Range("F14").Value = 1.988257407590117
$ws.Range("G14").Value = 0.03532932034805802
$ws.Range("E15").Value = 12.05397234475991
$ws.Range("F15").Value = 1.851495197540349
$ws.Range("G15").Value = 0.01628479628088075
$ws.Range("E16").Value = 11.65214072299821
$ws.Range("F16").Value = 1.777467942317213
$ws.Range("G16").Value = 0.04907795892936972
$ws.Range("E17").Value = 12.13036601629508
$ws.Range("F17").Value = 1.737891569528418
$ws.Range("G17").Value = 0.01005036882346722
$ws.Range("E18").Value = 11.84795243134278
$ws.Range("F18").Value = 2.023685607936933
$ws.Range("G18").Value = 0.03309791939920492
$ws.Range("E19").Value = 12.07687309199294
$ws.Range("F19").Value = 1.872127933119406
$ws.Range("G19").Value = 0.01441588430851781
$ws.Range("E20").Value = 11.63438580506858
$ws.Range("F20").Value = 1.78630491018048
$ws.Range("G20").Value = 0.05052692381900215
$ws.Range("E21").Value = 12.11002629600846
$ws.Range("F21").Value = 1.744151226340409
$ws.Range("G21").Value = 0.01171027740074582
$ws.Range("E22").Value = 11.81242656139587
$ws.Range("F22").Value = 2.088131162121978
$ws.Range("G22").Value = 0.0359971577077548
$ws.Range("E23").Value = 11.85753090958821
$ws.Range("F23").Value = 1.901330666879689
$ws.Range("G23").Value = 0.03231622732218487
$ws.Range("E24").Value = 11.50023580248815
$ws.Range("F24").Value = 1.796483169091391
$ws.Range("G24").Value = 0.06147480003299566
$ws.Range("E25").Value = 11.77939493382674
$ws.Range("F25").Value = 1.74175774948887
$ws.Range("G25").Value = 0.03869284285734509
$ws.Range("E26").Value = 11.8513352001668
$ws.Range("F26").Value = 2.139547219813363
$ws.Range("G26").Value = 0.03282185429571283
$ws.Range("E27").Value = 11.79501981542471
$ws.Range("F27").Value = 1.970909414250721
$ws.Range("G27").Value = 0.0374177085576618
$ws.Range("E28").Value = 11.47984161382105
$ws.Range("F28").Value = 1.809822873941755
$ws.Range("G28").Value = 0.06313915373197054
$ws.Range("E29").Value = 11.58564518082139
$ws.Range("F29").Value = 1.741694650585852
$ws.Range("G29").Value = 0.05450460783381283
